# Update the "Förändrad" (Changed) date column (C2:C6) from 2023-10-22 to 2023-10-25
# (Excel serial date 45221 -> 45224) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
